$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update betting odds values in row 2 as per the FlashScore data refresh
$ws.Range("G2").Value = 2.3
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 3.25
$ws.Range("K2").Value = 1.77
$ws.Range("L2").Value = 4.5
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("AH2").Value = 7
$ws.Range("AX2").Value = 23
$ws.Range("AZ2").Value = 81
